# Insert a new data row right after the existing row 34 (i.e. at row 35),
# shifting all subsequent rows (old 35..115) down to (new 36..116), and
# populate the newly inserted row 35 with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44791
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112035
$ws.Range("G35").Value = "Bruselas (repollito)"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 55
$ws.Range("K35").Value = 24000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = 24545
$ws.Range("N35").Value = "$/malla 10 kilos"
$ws.Range("O35").Value = "Provincia de Quillota"
$ws.Range("P35").Value = 2454
$ws.Range("Q35").Value = 10
$ws.Range("R35").Value = "Hortaliza"
